$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F, copying the header style from column E (bold, centered, bordered)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Update header: MSE -> RMSE, and set new column F header to "MSE_log"
$ws.Range("E1").Value = "RMSE"
$ws.Range("F1").Value = "MSE_log"

# Row 2 - Gradient Boosting
$ws.Range("B2").Value = 0.775719984354989
$ws.Range("C2").Value = 0.1784170197730449
$ws.Range("D2").Value = 0.806
$ws.Range("E2").Value = 2232917813.31779
$ws.Range("F2").Value = 0.0112439613306092

# Row 3 - Gradient Boosting (DeepWalk)
$ws.Range("B3").Value = 0.803511885674746
$ws.Range("C3").Value = 0.1689024077359783
$ws.Range("D3").Value = 0.79
$ws.Range("E3").Value = 2089996659.117547
$ws.Range("F3").Value = 0.009878683005223553

# Row 4 - Linear Regression
$ws.Range("B4").Value = 0.6447146954908073
$ws.Range("C4").Value = 0.3189993446933566
$ws.Range("D4").Value = 0.574
$ws.Range("E4").Value = 2810386636.278114
$ws.Range("F4").Value = 0.03115809645376205

# Row 5 - Linear Regression (DeepWalk)
$ws.Range("B5").Value = 0.6444440888307548
$ws.Range("C5").Value = 0.3190082859223518
$ws.Range("D5").Value = 0.5659999999999999
$ws.Range("E5").Value = 2811456712.318631
$ws.Range("F5").Value = 0.02986115581908628

# Row 6 - Random Forest
$ws.Range("B6").Value = 0.8181890743697067
$ws.Range("C6").Value = 0.1721062464955148
$ws.Range("D6").Value = 0.804
$ws.Range("E6").Value = 2010422973.744161
$ws.Range("F6").Value = 0.01033476672484333

# Row 7 - Random Forest (DeepWalk)
$ws.Range("B7").Value = 0.8090055490035011
$ws.Range("C7").Value = 0.1695272787264475
$ws.Range("D7").Value = 0.798
$ws.Range("E7").Value = 2060572143.494536
$ws.Range("F7").Value = 0.009944598743791801
